{"js": "// Update the title date and every divison-problem cell in the table.\nconst body = context.document.body;\n\n// --- 1. Title paragraph: \"2023-11-02 Thursday\" -> \"2023-11-03 Friday\" ---\nconst paras = body.paragraphs;\nparas.load(\"items/text\");\nawait context.sync();\n\nfor (const p of paras.items) {\n  if (p.text.indexOf(\"2023-11-02 Thursday\") !== -1) {\n    p.insertText(\"2023-11-03 Friday\", \"Replace\");\n    break;\n  }\n}\nawait context.sync();\n\n// --- 2. Table cells: replace each division problem by its row/col position ---\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Row -> [col0..col4] new values, matching the row indices that actually\n// contain text in the 20-row table (every 4th row starting at 0).\nconst rowValues = {\n  0: [\"24\u00f77=3, 3\", \"70\u00f75=14, 0\", \"47\u00f75=9, 2\", \"98\u00f75=19, 3\", \"44\u00f78=5, 4\"],\n  4: [\"96\u00f77=13, 5\", \"15\u00f73=5, 0\", \"26\u00f75=5, 1\", \"41\u00f77=5, 6\", \"48\u00f76=8, 0\"],\n  8: [\"16\u00f78=2, 0\", \"80\u00f73=26, 2\", \"14\u00f77=2, 0\", \"17\u00f72=8, 1\", \"28\u00f72=14, 0\"],\n  12: [\"53\u00f79=5, 8\", \"48\u00f77=6, 6\", \"54\u00f77=7, 5\", \"15\u00f74=3, 3\", \"58\u00f75=11, 3\"],\n  16: [\"76\u00f78=9, 4\", \"35\u00f79=3, 8\", \"64\u00f75=12, 4\", \"73\u00f78=9, 1\", \"87\u00f72=43, 1\"],\n};\n\nfor (const rowIndexStr of Object.keys(rowValues)) {\n  const rowIndex = parseInt(rowIndexStr, 10);\n  const values = rowValues[rowIndex];\n  for (let col = 0; col < values.length; col++) {\n    table.getCell(rowIndex, col).value = values[col];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the title date and every division-problem cell in the table.\n$d = $word.ActiveDocument\n\n# --- 1. Title paragraph: \"2023-11-02 Thursday\" -> \"2023-11-03 Friday\" ---\n$find = $d.Content.Find\n$find.Text = \"2023-11-02 Thursday\"\n$find.Replacement.Text = \"2023-11-03 Friday\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2) | Out-Null\n\n# --- 2. Table cells: replace each division problem by its row/col position ---\n# COM Cell(row, col) is 1-based; only every 4th row (1,5,9,13,17) holds data.\n$t = $d.Tables.Item(1)\n\n$t.Cell(1,1).Range.Text = \"24\u00f77=3, 3\"\n$t.Cell(1,2).Range.Text = \"70\u00f75=14, 0\"\n$t.Cell(1,3).Range.Text = \"47\u00f75=9, 2\"\n$t.Cell(1,4).Range.Text = \"98\u00f75=19, 3\"\n$t.Cell(1,5).Range.Text = \"44\u00f78=5, 4\"\n\n$t.Cell(5,1).Range.Text = \"96\u00f77=13, 5\"\n$t.Cell(5,2).Range.Text = \"15\u00f73=5, 0\"\n$t.Cell(5,3).Range.Text = \"26\u00f75=5, 1\"\n$t.Cell(5,4).Range.Text = \"41\u00f77=5, 6\"\n$t.Cell(5,5).Range.Text = \"48\u00f76=8, 0\"\n\n$t.Cell(9,1).Range.Text = \"16\u00f78=2, 0\"\n$t.Cell(9,2).Range.Text = \"80\u00f73=26, 2\"\n$t.Cell(9,3).Range.Text = \"14\u00f77=2, 0\"\n$t.Cell(9,4).Range.Text = \"17\u00f72=8, 1\"\n$t.Cell(9,5).Range.Text = \"28\u00f72=14, 0\"\n\n$t.Cell(13,1).Range.Text = \"53\u00f79=5, 8\"\n$t.Cell(13,2).Range.Text = \"48\u00f77=6, 6\"\n$t.Cell(13,3).Range.Text = \"54\u00f77=7, 5\"\n$t.Cell(13,4).Range.Text = \"15\u00f74=3, 3\"\n$t.Cell(13,5).Range.Text = \"58\u00f75=11, 3\"\n\n$t.Cell(17,1).Range.Text = \"76\u00f78=9, 4\"\n$t.Cell(17,2).Range.Text = \"35\u00f79=3, 8\"\n$t.Cell(17,3).Range.Text = \"64\u00f75=12, 4\"\n$t.Cell(17,4).Range.Text = \"73\u00f78=9, 1\"\n$t.Cell(17,5).Range.Text = \"87\u00f72=43, 1\"\n"}
